$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Roof Tile
$ws.Range("H19").Value = 1695.3334
$ws.Range("I19").Value = 1599.5
$ws.Range("J19").Value = 1743.25
$ws.Range("K19").Value = 1599.5
$ws.Range("L19").Value = 1743.25
$ws.Range("M19").Value = -1424.5
$ws.Range("N19").Value = -2093.25

# Row 51: Shark Oil
$ws.Range("H51").Value = 11332.833
$ws.Range("I51").Value = 10999.25
$ws.Range("K51").Value = 10999.25
$ws.Range("M51").Value = -10515.25

# Row 70: Holy Water
$ws.Range("H70").Value = 4769.647
$ws.Range("I70").Value = 3571.1538
$ws.Range("J70").Value = 5511.5713
$ws.Range("K70").Value = 10713.4614
$ws.Range("L70").Value = 16534.7139
$ws.Range("M70").Value = -10443.4614
$ws.Range("N70").Value = -17074.7139

# Row 73: Holy Water
$ws.Range("H73").Value = 4769.647
$ws.Range("I73").Value = 3571.1538
$ws.Range("J73").Value = 5511.5713
$ws.Range("K73").Value = 10713.4614
$ws.Range("L73").Value = 16534.7139
$ws.Range("M73").Value = -9777.4614
$ws.Range("N73").Value = -18406.7139

# Row 82: Draconian Potion of Dexterity
$ws.Range("H82").Value = 74000.336
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents() | Out-Null

# Row 85: Draconian Potion of Dexterity
$ws.Range("H85").Value = 74000.336
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents() | Out-Null

# Row 100: Beetle Glue
$ws.Range("H100").Value = 3334.25
$ws.Range("I100").Value = 3862.5
$ws.Range("K100").Value = 3862.5
$ws.Range("M100").Value = -3321.5

# Row 107: Enchanted Truegold Ink
$ws.Range("H107").Value = 1173.4546
$ws.Range("I107").Value = 863.5
$ws.Range("K107").Value = 863.5
$ws.Range("M107").Value = 1056.5

# Row 113: Starch Glue
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = 254

# Row 132: Growth Formula Lambda
$ws.Range("H132").Value = 6082.0835
$ws.Range("I132").Value = 6362.273
$ws.Range("K132").Value = 19086.819
$ws.Range("M132").Value = -16556.819

# Row 135: Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 694.8823
$ws.Range("I135").Value = 718.3125
$ws.Range("K135").Value = 6464.8125
$ws.Range("M135").Value = -3929.8125

$ws = $wb.Worksheets.Item("ARM")
# Row 11: Bronze Sollerets
$ws.Range("H11").Value = 15265
$ws.Range("J11").Value = 15265
$ws.Range("L11").Value = 15265
$ws.Range("N11").Value = -15553

# Row 35: Conical Alembic
$ws.Range("H35").Value = 2666.3333
$ws.Range("I35").Value = 2666.3333
$ws.Range("K35").Value = 2666.3333
$ws.Range("M35").Value = -2260.3333

# Row 122: High Durium Nugget
$ws.Range("H122").Value = 2214.2
$ws.Range("I122").Value = 2401.077
$ws.Range("K122").Value = 7203.231000000001
$ws.Range("M122").Value = -4753.231000000001

# Row 131: Chondrite Top of Maiming
$ws.Range("H131").Value = 71999.336
$ws.Range("J131").Value = 71999.336
$ws.Range("L131").Value = 71999.336
$ws.Range("N131").Value = -82079.336

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Iron Rivets
$ws.Range("H22").Value = 380.7
$ws.Range("I22").Value = 380.7
$ws.Range("K22").Value = 380.7
$ws.Range("M22").Value = -207.7

# Row 36: Iron Chocobotail Saw
$ws.Range("H36").Value = 2966.3333
$ws.Range("I36").Value = 2966.3333
$ws.Range("K36").Value = 2966.3333
$ws.Range("M36").Value = -2432.3333

# Row 86: Adamantite Nugget
$ws.Range("H86").Value = 8170.7144
$ws.Range("J86").Value = 14566.333
$ws.Range("L86").Value = 14566.333
$ws.Range("N86").Value = -16812.333

# Row 89: Adamantite Nugget
$ws.Range("H89").Value = 8170.7144
$ws.Range("J89").Value = 14566.333
$ws.Range("L89").Value = 72831.66500000001
$ws.Range("N89").Value = -84063.66500000001

# Row 95: High Steel Kris
$ws.Range("H95").Value = 14750
$ws.Range("J95").Value = 14750
$ws.Range("L95").Value = 14750
$ws.Range("N95").Value = -20242

# Row 105: Molybdenum Ingot
$ws.Range("H105").Value = 4239.7856
$ws.Range("I105").Value = 4034.4119
$ws.Range("K105").Value = 4034.4119
$ws.Range("M105").Value = -2287.4119

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Elm Lumber
$ws.Range("H22").Value = 3636799.5
$ws.Range("I22").Value = 556.6667
$ws.Range("K22").Value = 556.6667
$ws.Range("M22").Value = -206.6667

# Row 31: Walnut Lumber
$ws.Range("H31").Value = 5300
$ws.Range("I31").Value = 1920
$ws.Range("K31").Value = 1920
$ws.Range("M31").Value = -1625

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 5300
$ws.Range("I34").Value = 1920
$ws.Range("K34").Value = 1920
$ws.Range("M34").Value = -1718

# Row 99: Pine Lumber
$ws.Range("H99").Value = 2980.8
$ws.Range("I99").Value = 2980.8
$ws.Range("K99").Value = 2980.8
$ws.Range("M99").Value = -1482.8

# Row 108: White Oak Fishing Rod
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents() | Out-Null

# Row 116: Sandteak Rod
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents() | Out-Null

# Row 117: Sandteak Spinning Wheel
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").ClearContents() | Out-Null
$ws.Range("N117").ClearContents() | Out-Null

# Row 118: Sandteak Longbow
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents() | Out-Null

# Row 126: Red Pine Lumber
$ws.Range("H126").Value = 2980.8
$ws.Range("I126").Value = 2980.8
$ws.Range("K126").Value = 8942.400000000001
$ws.Range("M126").Value = -6472.400000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5: Maple Syrup
$ws.Range("H5").Value = 1211.5
$ws.Range("J5").Value = 1489.5
$ws.Range("L5").Value = 4468.5
$ws.Range("N5").Value = -4692.5

# Row 15: Grilled Carp
$ws.Range("H15").Value = 199.83333
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents() | Out-Null

# Row 51: Jerked Beef
$ws.Range("H51").Value = 796
$ws.Range("I51").Value = 796
$ws.Range("K51").Value = 2388
$ws.Range("M51").Value = -1928

# Row 135: Royal Maple Syrup
$ws.Range("H135").Value = 1211.5
$ws.Range("J135").Value = 1489.5
$ws.Range("L135").Value = 13405.5
$ws.Range("N135").Value = -18475.5

$ws = $wb.Worksheets.Item("GSM")
# Row 45: Mythril Ring
$ws.Range("H45").Value = 34999.25
$ws.Range("J45").Value = 34999.25
$ws.Range("L45").Value = 34999.25
$ws.Range("N45").Value = -36117.25

# Row 121: Petalite Bracelet of Fending
$ws.Range("H121").Value = 107555.445
$ws.Range("J121").Value = 107555.445
$ws.Range("L121").Value = 107555.445
$ws.Range("N121").Value = -111049.445

# Row 122: Ametrine
$ws.Range("H122").Value = 1822.619
$ws.Range("I122").Value = 1559.7646
$ws.Range("K122").Value = 4679.293799999999
$ws.Range("M122").Value = -2229.293799999999

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Aldgoat Leather
$ws.Range("H22").Value = 644.875
$ws.Range("I22").Value = 411.6875
$ws.Range("K22").Value = 411.6875
$ws.Range("M22").Value = -116.6875

# Row 27: Aldgoat Leather
$ws.Range("H27").Value = 644.875
$ws.Range("I27").Value = 411.6875
$ws.Range("K27").Value = 411.6875
$ws.Range("M27").Value = -304.6875

# Row 40: Toad Leather
$ws.Range("H40").Value = 4900
$ws.Range("J40").Value = 4900
$ws.Range("L40").Value = 4900
$ws.Range("N40").Value = -5172

# Row 46: Boar Leather
$ws.Range("H46").Value = 1941.1666
$ws.Range("I46").Value = 1474.1
$ws.Range("J46").Value = 4276.5
$ws.Range("K46").Value = 1474.1
$ws.Range("L46").Value = 4276.5
$ws.Range("M46").Value = -1286.1
$ws.Range("N46").Value = -4652.5

# Row 55: Peiste Leather
$ws.Range("H55").Value = 675.2353000000001
$ws.Range("I55").Value = 235.28572
$ws.Range("J55").Value = 983.2
$ws.Range("K55").Value = 235.28572
$ws.Range("L55").Value = 983.2
$ws.Range("M55").Value = -62.28572
$ws.Range("N55").Value = -1329.2

# Row 61: Raptor Leather
$ws.Range("H61").Value = 3996.5
$ws.Range("I61").Value = 3996.5
$ws.Range("K61").Value = 3996.5
$ws.Range("M61").Value = -3794.5

# Row 113: Atrociraptor Leather
$ws.Range("H113").Value = 3996.5
$ws.Range("I113").Value = 3996.5
$ws.Range("K113").Value = 3996.5
$ws.Range("M113").Value = -1826.5

# Row 130: Ophiotauroskin Boots of Healing
$ws.Range("H130").Value = 66661.664
$ws.Range("J130").Value = 66661.664
$ws.Range("L130").Value = 66661.664
$ws.Range("N130").Value = -76701.664

# Row 136: Br'aax Leather
$ws.Range("H136").Value = 2972.4285
$ws.Range("I136").Value = 2634.5
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7903.5
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5353.5
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
# Row 51: Linen Smock
$ws.Range("H51").Value = 21444
$ws.Range("I51").Value = 20666
$ws.Range("J51").Value = 23000
$ws.Range("K51").Value = 20666
$ws.Range("L51").Value = 23000
$ws.Range("M51").Value = -20156
$ws.Range("N51").Value = -24020

# Row 104: Twinsilk Turban of Aiming
$ws.Range("H104").Value = 24650
$ws.Range("J104").Value = 24650
$ws.Range("L104").Value = 24650
$ws.Range("N104").Value = -31638

# Row 130: AR-Caean Velvet Cap of Maiming
$ws.Range("H130").Value = 27999
$ws.Range("J130").Value = 27999
$ws.Range("L130").Value = 27999
$ws.Range("N130").Value = -38039

# Row 136: Sarcenet Cloth
$ws.Range("H136").Value = 5804.375
$ws.Range("I136").Value = 5359.625
$ws.Range("K136").Value = 16078.875
$ws.Range("M136").Value = -13528.875
